$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.550.77'
$ws.Range("E2").Value = '  +4.97%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.103.41'
$ws.Range("E3").Value = '  +3.02%  '
$ws.Range("E4").Value = '  +0.36%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.73'
$ws.Range("E5").Value = '  +3.67%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.99'
$ws.Range("E6").Value = '  +3.43%  '
$ws.Range("E7").Value = '  +0.24%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.092.10'
$ws.Range("E8").Value = '  +2.94%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.531'
$ws.Range("E9").Value = '  +1.55%  '
$ws.Range("E10").Value = '  +7.73%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.76'
$ws.Range("E11").Value = '  +9.48%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.468'
$ws.Range("E12").Value = '  +1.99%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000243'
$ws.Range("E13").Value = '  +4.62%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.54'
$ws.Range("E14").Value = '  +5.15%  '
$ws.Range("E15").Value = '  +0.79%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.619.30'
$ws.Range("E16").Value = '  +3.55%  '
$ws.Range("E17").Value = '  -0.49%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.105.69'
$ws.Range("E18").Value = '  +3.52%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '62.547.54'
$ws.Range("E19").Value = '  +5.59%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '453.77'
$ws.Range("E20").Value = '  +4.97%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.08'
$ws.Range("E21").Value = '  +2.41%  '
$ws.Range("E22").Value = '  +1.50%  '
$ws.Range("E23").Value = '  +5.50%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.70'
$ws.Range("E24").Value = '  +1.80%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '82.05'
$ws.Range("E25").Value = '  +1.46%  '
$ws.Range("E26").Value = '  +0.14%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.27'
$ws.Range("E27").Value = '  +2.77%  '
$ws.Range("E28").Value = '  +5.11%  '
$ws.Range("E29").Value = '  +0.65%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.23'
$ws.Range("E30").Value = '  +4.24%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.85'
$ws.Range("E31").Value = '  +12.13%  '
$ws.Range("E32").Value = '  +11.89%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.06'
$ws.Range("E33").Value = '  +4.45%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.05'
$ws.Range("E34").Value = '  +4.42%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0₃0800'
$ws.Range("E35").Value = '  +4.06%  '
$ws.Range("E36").Value = '  +1.58%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.20'
$ws.Range("E37").Value = '  +3.99%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '50.66'
$ws.Range("E38").Value = '  +3.44%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.02'
$ws.Range("E39").Value = '  +9.33%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.81'
$ws.Range("E40").Value = '  +1.35%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '423.50'
$ws.Range("E41").Value = '  +4.92%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.935.43'
$ws.Range("E42").Value = '  +5.94%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0372'
$ws.Range("E43").Value = '  +5.15%  '
$ws.Range("E44").Value = '  +11.10%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.111'
$ws.Range("E45").Value = '  +3.01%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.16'
$ws.Range("E46").Value = '  +7.12%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '125.44'
$ws.Range("E47").Value = '  +1.43%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '34.94'
$ws.Range("E49").Value = '  -2.69%  '
$ws.Range("E50").Value = '  +0.69%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '24.78'
$ws.Range("E51").Value = '  +4.97%  '
